$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "survey" sheet (sheet2.xml): add two new form sections ("visit" and
#    "plot"), mirroring the existing repeated block pattern (title row,
#    body/open-form row, exit-section row).
# ---------------------------------------------------------------------------
$wsSurvey = $wb.Worksheets.Item("survey")

# -- visit --
$wsSurvey.Rows.Item(53).RowHeight = 17.5
$wsSurvey.Cells.Item(53, 1).Value = "visit"

$wsSurvey.Rows.Item(54).RowHeight = 66
$wsSurvey.Cells.Item(54, 1).Style = $wsSurvey.Cells.Item(51, 1).Style
$wsSurvey.Cells.Item(54, 2).Value = "''?' + opendatakit.getHashString('../tables/visit/forms/visit/',null)"
$wsSurvey.Cells.Item(54, 5).Value = "external_link"
$wsSurvey.Cells.Item(54, 7).Value = "Open form"

$wsSurvey.Rows.Item(55).RowHeight = 17
$wsSurvey.Cells.Item(55, 1).Style = $wsSurvey.Cells.Item(52, 1).Style
$wsSurvey.Cells.Item(55, 2).Style = $wsSurvey.Cells.Item(52, 2).Style
$wsSurvey.Cells.Item(55, 3).Value = "exit section"

# -- plot --
$wsSurvey.Rows.Item(56).RowHeight = 17.5
$wsSurvey.Cells.Item(56, 1).Value = "plot"

$wsSurvey.Rows.Item(57).RowHeight = 66
$wsSurvey.Cells.Item(57, 1).Style = $wsSurvey.Cells.Item(54, 1).Style
$wsSurvey.Cells.Item(57, 2).Value = "''?' + opendatakit.getHashString('../tables/plot/forms/plot/',null)"
$wsSurvey.Cells.Item(57, 5).Value = "external_link"
$wsSurvey.Cells.Item(57, 7).Value = "Open form"

$wsSurvey.Rows.Item(58).RowHeight = 17
$wsSurvey.Cells.Item(58, 1).Style = $wsSurvey.Cells.Item(55, 1).Style
$wsSurvey.Cells.Item(58, 2).Style = $wsSurvey.Cells.Item(55, 2).Style
$wsSurvey.Cells.Item(58, 3).Value = "exit section"

# ---------------------------------------------------------------------------
# 2. "choices" sheet (sheet4.xml): register the two new test forms in the
#    "test_forms" choice list.
# ---------------------------------------------------------------------------
$wsChoices = $wb.Worksheets.Item("choices")

$wsChoices.Cells.Item(17, 1).Style = $wsChoices.Cells.Item(16, 1).Style
$wsChoices.Cells.Item(17, 1).Value = "test_forms"
$wsChoices.Cells.Item(17, 2).Value = "visit"
$wsChoices.Cells.Item(17, 3).Value = "Visit"

$wsChoices.Cells.Item(18, 1).Style = $wsChoices.Cells.Item(17, 1).Style
$wsChoices.Cells.Item(18, 1).Value = "test_forms"
$wsChoices.Cells.Item(18, 2).Value = "plot"
$wsChoices.Cells.Item(18, 3).Value = "Plot"

# ---------------------------------------------------------------------------
# 3. Update view/selection state to reflect the new bottom-of-sheet rows.
#    ("choices" stays the active/tabSelected sheet throughout, matching the
#    before/after workbook state.)
# ---------------------------------------------------------------------------
$wsSurvey.Activate()
$wsSurvey.Range("B58").Select()

$wsChoices.Activate()
$wsChoices.Range("B19").Select()

Write-Output "done"
